$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "O8"  = 1.25
    "P8"  = 3.75
    "Q8"  = 1.9
    "R8"  = 1.95

    "G9"  = 2.35
    "I9"  = 2.75
    "AA9" = 17
    "AN9" = 4.5
    "AZ9" = 51
    "BD9" = 151

    "Q13" = 1.53
    "R13" = 2.4

    "G17"  = 2.63
    "H17"  = 3.2
    "I17"  = 2.45
    "J17"  = 3.4
    "L17"  = 3.25
    "M17"  = 1.06
    "N17"  = 10
    "Q17"  = 2.1
    "R17"  = 1.7
    "AG17" = 8
    "AJ17" = 26
    "AO17" = 15
    "AY17" = 26

    "G18"  = 2.62
    "H18"  = 3.1
    "I18"  = 2.65
    "J18"  = 3.1
    "L18"  = 3.15
    "O18"  = 1.34
    "P18"  = 3
    "Q18"  = 2.02
    "R18"  = 1.75
    "T18"  = 2.72
    "U18"  = 1.8
    "V18"  = 1.9
    "W18"  = 8
    "Y18"  = 9.75
    "AB18" = 32
    "AE18" = 15
    "AF18" = 70
    "AG18" = 8
    "AH18" = 12.5
    "AI18" = 10
    "AJ18" = 28
    "AM18" = 600
    "AP18" = 22
    "AR18" = 90
    "AT18" = 2.72
    "AU18" = 7.2
    "AV18" = 65
    "AW18" = 4.5

    "N20" = 8

    "AT22" = 2.63

    "AT24" = 2.63

    "Q26" = 1.36
    "R26" = 3.1

    "G28"  = 2.37
    "H28"  = 3.2
    "I28"  = 2.72
    "J28"  = 3.05
    "L28"  = 3.35
    "W28"  = 8.25
    "X28"  = 12
    "Y28"  = 9.25
    "Z28"  = 26
    "AA28" = 20
    "AB28" = 28
    "AG28" = 9.25
    "AH28" = 14.5
    "AI28" = 10
    "AJ28" = 32
    "AK28" = 23
    "AM28" = 450
    "AN28" = 4.4
    "AO28" = 13
    "AP28" = 21
    "AQ28" = 55
    "AR28" = 90
    "AS28" = 250
    "AW28" = 4.75
    "AX28" = 15
    "AZ28" = 70
    "BB28" = 250
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
